$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Final content of the "Milepæle" table (A:B existing columns, C:D new) ---
# Row 1 is the header row; rows 2-8 are data rows.
# Three brand-new milestone rows ("Interview guide", "Udførelse af interview",
# "Interview") are inserted above the existing "Problem analyse" row, and the
# wording of the "Problem analyse" deadline cell is updated.

$ws.Range("A2").Value = "Interview guide"
$ws.Range("B2").Value = "Laver interview guide onsdag d. 15/10"

$ws.Range("A3").Value = "Udførelse af interview"
$ws.Range("B3").Value = "Så hurtigt som muligt, gerne inden d. 20/10"

$ws.Range("A4").Value = "Interview"
$ws.Range("B4").Value = "Sendes til Mona d. 20/10 kl. 12.00"

$ws.Range("A5").Value = "Problem analyse"
$ws.Range("B5").Value = "Problem analysen SKAL sendes til Mona 30/12 kl. 12"

$ws.Range("A6").Value = "Påbegyndelse af programmering"
$ws.Range("B6").Value = "Programmeringen skal gerne være påbegyndt inden statusseminariet"

$ws.Range("A7").Value = "Rapport klar til at rettes"
$ws.Range("B7").Value = "Skal være klar til at blive rettet mandag d. 8/12 eller onsdag 10/12"

$ws.Range("A8").Value = "Rapport aflevering"
$ws.Range("B8").Value = "Rapporten skal afleveres torsdag d. 18/12 kl. 12.00"

# --- Resize the structured table to cover the new rows/columns first, so
# that writing the new header text below correctly syncs into table1.xml
# (renaming ListColumns after the fact does not update the saved table part).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D8"))

# New columns for tracking interview progress.
$ws.Range("C1").Value = "Klaret/ikke klaret"
$ws.Range("D1").Value = "Hvornår blev det færdigt?"

# Bold the milestone-name column for every data row except the very last one
# (matches the look of the original "Milepæle" column styling).
$ws.Range("A2:A7").Font.Bold = $true
$ws.Range("A2:A7").Font.Name = "Calibri"
$ws.Range("A2:A7").Font.Size = 11
$ws.Range("A8").Font.Bold = $true

# --- Column widths for the two new columns ---
$ws.Columns.Item(3).ColumnWidth = 15.4
$ws.Columns.Item(4).ColumnWidth = 21.6

# --- Selection matches the author's last active cell ---
$ws.Range("C3").Select() | Out-Null
